# feat: add DPYD variant
# Appends the new DPYD haplotype row to the "manual_snps" sheet, mirroring
# the manual entry of a new rs_id / gene / notes record at the bottom of
# the existing table (row 30, right under the previous last row, 29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manual_snps")

$ws.Range("A30").Value = "rs75017182"
$ws.Range("B30").Value = "DPYD"
$ws.Range("C30").Value = "Haplotype B3: c.1236G>A; c.1129-5923C>G - Compound variant"

# Match the formatting used by the rest of the data rows (Arial 12, not bold).
$rng = $ws.Range("A30:C30")
$rng.Font.Name = "Arial"
$rng.Font.Size = 12
$rng.Font.Bold = $false
